# Update the chip mapping table (DBInfo sheet) with new DB build numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBInfo")

# Merlin8 (row 3) now has an mpbackup entry: DB2919.
$ws.Range("P3").Value = "DB2919"
$ws.Range("R3").Value = "/DailyBuild/Merlin8/DB2919_Merlin8_64Bit_Android14_Ref_Plus_Wave_Backup_GoogleGMS"
$ws.Range("Q3").Value = "DB2919_Merlin8_64Bit_Android14_Ref_Plus_Wave_Backup_GoogleGMS"

# mp_DB_Info (column L) gets populated for the other chips.
$ws.Range("L2").Value = "DB2589"   # Merlin7
$ws.Range("L4").Value = "DB2897"   # Merlin9
$ws.Range("L5").Value = "DB2588"   # Mac7p
$ws.Range("L6").Value = "DB2592"   # Mac8q

# The new, much longer mpbackup folder/path strings no longer fit in the
# old column widths, so widen columns Q (mpbackup_DB_Folder) and R
# (mpbackup_SftpPath) to fit the content.
$ws.Range("Q1").EntireColumn.ColumnWidth = 67.14285714285714
$ws.Range("R1").EntireColumn.ColumnWidth = 85.71428571428571

# Reset the view: scroll back to the left and move the selection.
$ws.Range("E10").Select()

# Configure the print setup (A4, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
